$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 46.666668
$ws.Range("I2").Value = 44.333332
$ws.Range("J2").Value = 49
$ws.Range("K2").Value = 44.333332
$ws.Range("L2").Value = 49
$ws.Range("M2").Value = 68.666668
$ws.Range("N2").Value = -275
$ws.Range("H18").Value = 1290.5454
$ws.Range("I18").Value = 1119.6
$ws.Range("J18").Value = 3000
$ws.Range("K18").Value = 1119.6
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -835.5999999999999
$ws.Range("N18").Value = -3568
$ws.Range("H32").Value = 148623.88
$ws.Range("I32").Value = 17497.5
$ws.Range("J32").Value = 192332.67
$ws.Range("K32").Value = 17497.5
$ws.Range("L32").Value = 192332.67
$ws.Range("M32").Value = -17171.5
$ws.Range("N32").Value = -192984.67
$ws.Range("H33").Value = 364
$ws.Range("I33").Value = 276.6154
$ws.Range("J33").Value = 1500
$ws.Range("K33").Value = 276.6154
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = -47.61540000000002
$ws.Range("N33").Value = -1958
$ws.Range("H76").Value = 125003200
$ws.Range("I76").Value = 166669600
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 166669600
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -166669285
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 125003200
$ws.Range("I79").Value = 166669600
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 166669600
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -166668508
$ws.Range("N79").Value = -6184
$ws.Range("H103").Value = 187.14285
$ws.Range("I103").Value = 176.83333
$ws.Range("J103").Value = 249
$ws.Range("K103").Value = 530.49999
$ws.Range("L103").Value = 747
$ws.Range("M103").Value = 55.50000999999997
$ws.Range("N103").Value = -1919
$ws.Range("H116").Value = 3037109.2
$ws.Range("I116").Value = 6136.857
$ws.Range("J116").Value = 8341311
$ws.Range("K116").Value = 6136.857
$ws.Range("L116").Value = 8341311
$ws.Range("M116").Value = -2694.857
$ws.Range("N116").Value = -8348195
$ws.Range("H133").Value = 91459.375
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 91459.375
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 91459.375
$ws.Range("N133").Value = -101579.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11374937
$ws.Range("I45").Value = 5718.5
$ws.Range("J45").Value = 25018000
$ws.Range("K45").Value = 5718.5
$ws.Range("L45").Value = 25018000
$ws.Range("M45").Value = -5341.5
$ws.Range("N45").Value = -25018754
$ws.Range("H97").Value = 607.17645
$ws.Range("I97").Value = 457.625
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 457.625
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = 38.375
$ws.Range("N97").Value = -3992
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2594.25
$ws.Range("I132").Value = 1839.3846
$ws.Range("J132").Value = 3996.1428
$ws.Range("K132").Value = 5518.1538
$ws.Range("L132").Value = 11988.4284
$ws.Range("M132").Value = -2988.1538
$ws.Range("N132").Value = -17048.4284
$ws.Range("H139").Value = 62076
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 62076
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 62076
$ws.Range("N139").Value = -72356
$ws.Range("H141").Value = 89000
$ws.Range("I141").Value = 88000
$ws.Range("J141").Value = 89500
$ws.Range("K141").Value = 88000
$ws.Range("L141").Value = 89500
$ws.Range("M141").Value = -82820
$ws.Range("N141").Value = -99860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1943.0588
$ws.Range("I134").Value = 1136.25
$ws.Range("J134").Value = 3879.4
$ws.Range("K134").Value = 3408.75
$ws.Range("L134").Value = 11638.2
$ws.Range("M134").Value = -873.75
$ws.Range("N134").Value = -16708.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 749.8570999999999
$ws.Range("I22").Value = 199.66667
$ws.Range("J22").Value = 1162.5
$ws.Range("K22").Value = 199.66667
$ws.Range("L22").Value = 1162.5
$ws.Range("M22").Value = 150.33333
$ws.Range("N22").Value = -1862.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 9997.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 9997.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 29992.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -31364.5
$ws.Range("H65").Value = 9997.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 9997.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 89977.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -96841.5
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H97").Value = 127.55556
$ws.Range("I97").Value = 112.5
$ws.Range("J97").Value = 139.6
$ws.Range("K97").Value = 337.5
$ws.Range("L97").Value = 418.8
$ws.Range("M97").Value = 158.5
$ws.Range("N97").Value = -1410.8
$ws.Range("H112").Value = 5061.931
$ws.Range("I112").Value = 2999.3333
$ws.Range("H117").Value = 1194.8
$ws.Range("I117").Value = 761.2
$ws.Range("J117").Value = 1628.4
$ws.Range("K117").Value = 2283.6
$ws.Range("L117").Value = 4885.200000000001
$ws.Range("M117").Value = 1158.4
$ws.Range("N117").Value = -11769.2
$ws.Range("H121").Value = 2386.125
$ws.Range("I121").Value = 1162.5
$ws.Range("J121").Value = 2560.9285
$ws.Range("K121").Value = 3487.5
$ws.Range("L121").Value = 7682.7855
$ws.Range("M121").Value = -2177.5
$ws.Range("N121").Value = -10302.7855
$ws.Range("H131").Value = 1219.8125
$ws.Range("I131").Value = 831.63635
$ws.Range("J131").Value = 2073.8
$ws.Range("K131").Value = 2494.90905
$ws.Range("L131").Value = 6221.400000000001
$ws.Range("M131").Value = 2545.09095
$ws.Range("N131").Value = -16301.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 11291.833
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 11291.833
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 11291.833
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -15035.833
$ws.Range("H99").Value = 13042.923
$ws.Range("I99").Value = 6819.75
$ws.Range("J99").Value = 23000
$ws.Range("K99").Value = 6819.75
$ws.Range("L99").Value = 23000
$ws.Range("M99").Value = -4573.75
$ws.Range("N99").Value = -27492
$ws.Range("H109").Value = 70251.336
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 70251.336
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 70251.336
$ws.Range("N109").Value = -72331.336
$ws.Range("H132").Value = 5751.5625
$ws.Range("I132").Value = 4826.375
$ws.Range("J132").Value = 6676.75
$ws.Range("K132").Value = 14479.125
$ws.Range("L132").Value = 20030.25
$ws.Range("M132").Value = -11949.125
$ws.Range("N132").Value = -25090.25
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7271.6875
$ws.Range("I22").Value = 1299.5714
$ws.Range("J22").Value = 11916.667
$ws.Range("K22").Value = 1299.5714
$ws.Range("L22").Value = 11916.667
$ws.Range("M22").Value = -1004.5714
$ws.Range("N22").Value = -12506.667
$ws.Range("H27").Value = 7271.6875
$ws.Range("I27").Value = 1299.5714
$ws.Range("J27").Value = 11916.667
$ws.Range("K27").Value = 1299.5714
$ws.Range("L27").Value = 11916.667
$ws.Range("M27").Value = -1192.5714
$ws.Range("N27").Value = -12130.667
$ws.Range("H46").Value = 10538.385
$ws.Range("I46").Value = 16142.857
$ws.Range("J46").Value = 3999.8333
$ws.Range("K46").Value = 16142.857
$ws.Range("L46").Value = 3999.8333
$ws.Range("M46").Value = -15954.857
$ws.Range("N46").Value = -4375.8333
$ws.Range("H134").Value = 139638.67
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 139638.67
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 139638.67
$ws.Range("N134").Value = -149778.67
$ws.Range("H135").Value = 80197
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 80197
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 80197
$ws.Range("N135").Value = -90337
$ws.Range("H138").Value = 140760.28
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 140760.28
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 140760.28
$ws.Range("N138").Value = -151040.28
$ws.Range("H141").Value = 132500
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 132500
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 132500
$ws.Range("N141").Value = -142860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 206010060
$ws.Range("I2").Value = 206010060
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 206010060
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -206009948
$ws.Range("H4").Value = 19555
$ws.Range("I4").Value = 20888.666
$ws.Range("J4").Value = 18554.75
$ws.Range("K4").Value = 20888.666
$ws.Range("L4").Value = 18554.75
$ws.Range("M4").Value = -20775.666
$ws.Range("N4").Value = -18780.75
$ws.Range("H46").Value = 99884.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 99884.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 99884.5
$ws.Range("N46").Value = -100346.5
$ws.Range("H126").Value = 3559.25
$ws.Range("I126").Value = 1993.5
$ws.Range("J126").Value = 5125
$ws.Range("K126").Value = 5980.5
$ws.Range("L126").Value = 15375
$ws.Range("M126").Value = -3510.5
$ws.Range("N126").Value = -20315
$ws.Range("H133").Value = 70335.25
$ws.Range("I133").Value = 75500
$ws.Range("J133").Value = 68613.664
$ws.Range("K133").Value = 75500
$ws.Range("L133").Value = 68613.664
$ws.Range("M133").Value = -70440
$ws.Range("N133").Value = -78733.664
$ws.Range("H134").Value = 99884.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 99884.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 299653.5
$ws.Range("N134").Value = -304723.5
$ws.Range("H135").Value = 61656.125
$ws.Range("I135").Value = 59500
$ws.Range("J135").Value = 62374.832
$ws.Range("K135").Value = 59500
$ws.Range("L135").Value = 62374.832
$ws.Range("M135").Value = -54430
$ws.Range("N135").Value = -72514.83199999999
$ws.Range("H137").Value = 149333
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 149333
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 149333
$ws.Range("N137").Value = -159533
$ws.Range("H140").Value = 150000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 150000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 150000
$ws.Range("N140").Value = -160360
$ws.Range("H141").Value = 62926.875
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 62926.875
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 62926.875
$ws.Range("N141").Value = -73286.875

Write-Host "Updated cells: 307 set, 7 cleared"
